# Auto-generated Excel COM-interop script
# Applies cryptocurrency price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.561.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.944.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.17'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.70%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0809'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.10'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.229.18'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.813'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.42'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.20'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.938.73'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.538.76'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.37'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0857'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.04'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.99'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.30'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.08'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.135'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +16.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.25'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.119'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.67'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.18'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.24'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.69%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.32'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +16.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0993'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.69%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.15'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.91'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.25%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.03'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.342.61'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.72'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.14'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.38%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.122.01'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.12'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.15%  '
